# Apply updated odds values (and the shifted Odd_CS_* header labels)
# from the 2025-03-13 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels - Odd_CS_4-4 moved from AO1 to AI1, shifting AI1:AN1 -> AJ1:AO1
$ws.Range("AI1").Value = "Odd_CS_4-4"
$ws.Range("AJ1").Value = "Odd_CS_0-1"
$ws.Range("AK1").Value = "Odd_CS_0-2"
$ws.Range("AL1").Value = "Odd_CS_1-2"
$ws.Range("AM1").Value = "Odd_CS_0-3"
$ws.Range("AN1").Value = "Odd_CS_1-3"
$ws.Range("AO1").Value = "Odd_CS_2-3"

# Row 2: ITMjSKgS (El Gaish vs El Ismaily) odds refresh
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 2.77
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 2.95
$ws.Range("K2").Value = 1.87
$ws.Range("L2").Value = 4.25
$ws.Range("P2").Value = 2.32
$ws.Range("Q2").Value = 2.57
$ws.Range("R2").Value = 1.45
$ws.Range("S2").Value = 4.6
$ws.Range("T2").Value = 1.16
$ws.Range("W2").Value = 2.1
$ws.Range("X2").Value = 1.65
$ws.Range("Y2").Value = 5.9
$ws.Range("Z2").Value = 10
$ws.Range("AB2").Value = 24
$ws.Range("AD2").Value = 40
$ws.Range("AF2").Value = 5.6
$ws.Range("AG2").Value = 17.5
$ws.Range("AH2").Value = 110
$ws.Range("AI2").Value = 900
$ws.Range("AJ2").Value = 7.3
$ws.Range("AK2").Value = 16.5
$ws.Range("AL2").Value = 13
$ws.Range("AM2").Value = 55
$ws.Range("AN2").Value = 45
$ws.Range("AO2").Value = 60

# Row 3: 6NrDPJR1 (Arema FC vs Barito Putera) odds refresh
$ws.Range("G3").Value = 1.83
$ws.Range("I3").Value = 3.55
$ws.Range("J3").Value = 2.35
$ws.Range("L3").Value = 3.85
$ws.Range("Q3").Value = 1.6
$ws.Range("R3").Value = 2.07
$ws.Range("S3").Value = 2.4
$ws.Range("W3").Value = 1.57
$ws.Range("X3").Value = 2.12
$ws.Range("Y3").Value = 9.25
$ws.Range("Z3").Value = 10
$ws.Range("AB3").Value = 16
$ws.Range("AE3").Value = 14
$ws.Range("AG3").Value = 13.5
$ws.Range("AH3").Value = 50
$ws.Range("AI3").Value = 300
$ws.Range("AJ3").Value = 13.5
$ws.Range("AK3").Value = 21
$ws.Range("AL3").Value = 12
$ws.Range("AM3").Value = 50
$ws.Range("AN3").Value = 28
$ws.Range("AO3").Value = 30

# Row 4: GQ9I3mnI (Al Fateh vs Al Raed) odds refresh
$ws.Range("G4").Value = 1.62
$ws.Range("H4").Value = 4.1
$ws.Range("I4").Value = 4.75
$ws.Range("J4").Value = 2.1
$ws.Range("K4").Value = 2.38
$ws.Range("L4").Value = 4.75
$ws.Range("N4").Value = 15
$ws.Range("Q4").Value = 1.62
$ws.Range("R4").Value = 2.25
$ws.Range("S4").Value = 2.5
$ws.Range("T4").Value = 1.5
$ws.Range("U4").Value = 1.3
$ws.Range("V4").Value = 3.4
$ws.Range("AD4").Value = 21
$ws.Range("AE4").Value = 15
$ws.Range("AI4").Value = 151
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 26
$ws.Range("AL4").Value = 15
$ws.Range("AM4").Value = 51
$ws.Range("AO4").Value = 34
$ws.Range("AP4").Value = 2.03
$ws.Range("AQ4").Value = 1.78

# Row 5: xKMjklXO (Al Ittihad vs Al Riyadh) odds refresh
$ws.Range("G5").Value = 1.4
$ws.Range("I5").Value = 6.25
$ws.Range("L5").Value = 5.5
$ws.Range("M5").Value = 21
$ws.Range("N5").Value = 1.03
$ws.Range("Q5").Value = 1.4
$ws.Range("AB5").Value = 11
$ws.Range("AI5").Value = 126
$ws.Range("AJ5").Value = 23
$ws.Range("AK5").Value = 41
$ws.Range("AL5").Value = 19
$ws.Range("AM5").Value = 67
$ws.Range("AO5").Value = 41

# Row 6: lbj97Vng (Al Shabab vs Al Orubah) odds refresh
$ws.Range("G6").Value = 1.25
$ws.Range("H6").Value = 6.25
$ws.Range("I6").Value = 9.5
$ws.Range("J6").Value = 1.62
$ws.Range("K6").Value = 2.75
$ws.Range("L6").Value = 8
$ws.Range("Q6").Value = 1.44
$ws.Range("R6").Value = 2.63
$ws.Range("W6").Value = 1.91
$ws.Range("X6").Value = 1.8
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 7
$ws.Range("AA6").Value = 10
$ws.Range("AB6").Value = 8
$ws.Range("AD6").Value = 26
$ws.Range("AG6").Value = 23
$ws.Range("AI6").Value = 600
$ws.Range("AJ6").Value = 26
$ws.Range("AK6").Value = 41
$ws.Range("AL6").Value = 26
$ws.Range("AM6").Value = 101
$ws.Range("AO6").Value = 51

# Row 7: vkY2Kc5m (Nafta vs Domzale) odds refresh
$ws.Range("G7").Value = 1.9
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 2.42
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 4.35
$ws.Range("N7").Value = 6.6
$ws.Range("O7").Value = 1.35
$ws.Range("P7").Value = 2.95
$ws.Range("Q7").Value = 2.02
$ws.Range("R7").Value = 1.7
$ws.Range("S7").Value = 3.4
$ws.Range("T7").Value = 1.27
$ws.Range("U7").Value = 1.38
$ws.Range("V7").Value = 2.82
$ws.Range("W7").Value = 1.85
$ws.Range("X7").Value = 1.85
$ws.Range("Y7").Value = 6.5
$ws.Range("Z7").Value = 8.5
$ws.Range("AB7").Value = 16.5
$ws.Range("AD7").Value = 29
$ws.Range("AE7").Value = 6.6
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 75
$ws.Range("AI7").Value = 600
$ws.Range("AJ7").Value = 10.75
$ws.Range("AK7").Value = 22
$ws.Range("AL7").Value = 13.5
$ws.Range("AM7").Value = 65
$ws.Range("AN7").Value = 40
$ws.Range("AO7").Value = 45
